$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 547.8823
$ws.Range("I33").Value = 635.53845
$ws.Range("J33").Value = 263
$ws.Range("K33").Value = 635.53845
$ws.Range("L33").Value = 263
$ws.Range("M33").Value = -406.53845
$ws.Range("N33").Value = -721
$ws.Range("H137").Value = 2138577.5
$ws.Range("I137").Value = 4387502
$ws.Range("J137").Value = 2099.1
$ws.Range("K137").Value = 13162506
$ws.Range("L137").Value = 6297.299999999999
$ws.Range("M137").Value = -13159956
$ws.Range("N137").Value = -11397.3
$ws.Range("H138").Value = 4156.9277
$ws.Range("J138").Value = 4230.523
$ws.Range("L138").Value = 12691.569
$ws.Range("N138").Value = -22971.569
$ws.Range("H141").Value = 3019.2222
$ws.Range("I141").Value = 1411.1428
$ws.Range("J141").Value = 8647.5
$ws.Range("K141").Value = 4233.428400000001
$ws.Range("L141").Value = 25942.5
$ws.Range("M141").Value = 946.5715999999993
$ws.Range("N141").Value = -36302.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15876487
$ws.Range("I61").Value = 33335472
$ws.Range("J61").Value = 4683.091
$ws.Range("K61").Value = 33335472
$ws.Range("L61").Value = 4683.091
$ws.Range("M61").Value = -33335260
$ws.Range("N61").Value = -5107.091
$ws.Range("H62").Value = 32749.666
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 32749.666
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H74").Value = 14287801
$ws.Range("I74").Value = 1172.45
$ws.Range("J74").Value = 33336638
$ws.Range("K74").Value = 1172.45
$ws.Range("L74").Value = 33336638
$ws.Range("M74").Value = -298.45
$ws.Range("N74").Value = -33338386
$ws.Range("H77").Value = 14287801
$ws.Range("I77").Value = 1172.45
$ws.Range("J77").Value = 33336638
$ws.Range("K77").Value = 5862.25
$ws.Range("L77").Value = 166683190
$ws.Range("M77").Value = -1494.25
$ws.Range("N77").Value = -166691926
$ws.Range("H122").Value = 73099.42999999999
$ws.Range("I122").Value = 73099.42999999999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 219298.29
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -216848.29
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 34714.5
$ws.Range("J123").Value = 34714.5
$ws.Range("L123").Value = 34714.5
$ws.Range("N123").Value = -44514.5
$ws.Range("H136").Value = 15876487
$ws.Range("I136").Value = 33335472
$ws.Range("J136").Value = 4683.091
$ws.Range("K136").Value = 100006416
$ws.Range("L136").Value = 14049.273
$ws.Range("M136").Value = -100003866
$ws.Range("N136").Value = -19149.273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3245.0857
$ws.Range("I134").Value = 3328.3076
$ws.Range("K134").Value = 9984.9228
$ws.Range("M134").Value = -7449.9228

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6206.5913
$ws.Range("I31").Value = 2639.2188
$ws.Range("J31").Value = 8078
$ws.Range("K31").Value = 2639.2188
$ws.Range("L31").Value = 8078
$ws.Range("M31").Value = -2344.2188
$ws.Range("N31").Value = -8668
$ws.Range("H34").Value = 6206.5913
$ws.Range("I34").Value = 2639.2188
$ws.Range("J34").Value = 8078
$ws.Range("K34").Value = 2639.2188
$ws.Range("L34").Value = 8078
$ws.Range("M34").Value = -2437.2188
$ws.Range("N34").Value = -8482
$ws.Range("H58").Value = 2226.25
$ws.Range("I58").Value = 2164.182
$ws.Range("K58").Value = 2164.182
$ws.Range("M58").Value = -1961.182
$ws.Range("H132").Value = 17545986
$ws.Range("I132").Value = 31252074
$ws.Range("J132").Value = 7577922
$ws.Range("K132").Value = 93756222
$ws.Range("L132").Value = 22733766
$ws.Range("M132").Value = -93753692
$ws.Range("N132").Value = -22738826
$ws.Range("H134").Value = 8626840
$ws.Range("I134").Value = 12506936
$ws.Range("K134").Value = 37520808
$ws.Range("M134").Value = -37518273
$ws.Range("H136").Value = 2226.25
$ws.Range("I136").Value = 2164.182
$ws.Range("K136").Value = 6492.545999999999
$ws.Range("M136").Value = -3942.545999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1243.9535
$ws.Range("J5").Value = 1633.1177
$ws.Range("L5").Value = 4899.3531
$ws.Range("N5").Value = -5123.3531
$ws.Range("H113").Value = 539.38806
$ws.Range("I113").Value = 509.17776
$ws.Range("J113").Value = 601.1818
$ws.Range("K113").Value = 1527.53328
$ws.Range("L113").Value = 1803.5454
$ws.Range("M113").Value = 642.4667200000001
$ws.Range("N113").Value = -6143.5454
$ws.Range("H117").Value = 2660.24
$ws.Range("I117").Value = 1892.875
$ws.Range("J117").Value = 3021.353
$ws.Range("K117").Value = 5678.625
$ws.Range("L117").Value = 9064.059000000001
$ws.Range("M117").Value = -2236.625
$ws.Range("N117").Value = -15948.059
$ws.Range("H129").Value = 820372.25
$ws.Range("I129").Value = 521.55554
$ws.Range("J129").Value = 1083895.6
$ws.Range("K129").Value = 1564.66662
$ws.Range("L129").Value = 3251686.8
$ws.Range("M129").Value = 3435.33338
$ws.Range("N129").Value = -3261686.8
$ws.Range("H135").Value = 1243.9535
$ws.Range("J135").Value = 1633.1177
$ws.Range("L135").Value = 14698.0593
$ws.Range("N135").Value = -19768.0593
$ws.Range("H139").Value = 273416.94
$ws.Range("I139").Value = 358704.28
$ws.Range("J139").Value = 8078.5557
$ws.Range("K139").Value = 1076112.84
$ws.Range("L139").Value = 24235.6671
$ws.Range("M139").Value = -1070972.84
$ws.Range("N139").Value = -34515.6671
$ws.Range("H140").Value = 1431.7069
$ws.Range("I140").Value = 945.65
$ws.Range("J140").Value = 2511.8333
$ws.Range("K140").Value = 2836.95
$ws.Range("L140").Value = 7535.499899999999
$ws.Range("M140").Value = 2343.05
$ws.Range("N140").Value = -17895.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 29933.334
$ws.Range("J62").Value = 29933.334
$ws.Range("L62").Value = 29933.334
$ws.Range("N62").Value = -31305.334
$ws.Range("H65").Value = 29933.334
$ws.Range("J65").Value = 29933.334
$ws.Range("L65").Value = 89800.00199999999
$ws.Range("N65").Value = -96664.00199999999
$ws.Range("H122").Value = 2975
$ws.Range("I122").Value = 3133.3333
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 9399.999899999999
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -6949.999899999999
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 33339402
$ws.Range("I132").Value = 55563116
$ws.Range("K132").Value = 166689348
$ws.Range("M132").Value = -166686818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 22857.143
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 22857.143
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 22857.143
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -32657.143
$ws.Range("H132").Value = 2752.6956
$ws.Range("I132").Value = 2021.32
$ws.Range("J132").Value = 3623.3809
$ws.Range("K132").Value = 6063.96
$ws.Range("L132").Value = 10870.1427
$ws.Range("M132").Value = -3533.96
$ws.Range("N132").Value = -15930.1427
$ws.Range("H140").Value = 68610.625
$ws.Range("J140").Value = 68610.625
$ws.Range("L140").Value = 68610.625
$ws.Range("N140").Value = -78970.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 57724.332
$ws.Range("J76").Value = 57724.332
$ws.Range("L76").Value = 57724.332
$ws.Range("N76").Value = -58354.332
$ws.Range("H79").Value = 57724.332
$ws.Range("J79").Value = 57724.332
$ws.Range("L79").Value = 57724.332
$ws.Range("N79").Value = -59908.332
$ws.Range("H132").Value = 4631848.5
$ws.Range("I132").Value = 2115.72
$ws.Range("K132").Value = 6347.16
$ws.Range("M132").Value = -3817.16
$ws.Range("H136").Value = 4590.773
$ws.Range("I136").Value = 6868.8
$ws.Range("J136").Value = 2692.4167
$ws.Range("K136").Value = 20606.4
$ws.Range("L136").Value = 8077.250100000001
$ws.Range("M136").Value = -18056.4
$ws.Range("N136").Value = -13177.2501
